# Update the cryptos list worksheet with the latest scraped values.
# Prices (column D) and volume deltas (column E) are stored as plain text
# in the workbook (not numbers), so every write goes through a small
# helper that forces the cell to remain text - even when the new value
# looks like a number (e.g. "1.00" or "418.97") - while restoring the
# cell's original ("Normal") style afterwards so no stray number format
# is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "71.220.05"
Set-TextValue "E2" "  +0.97%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.849.26"
Set-TextValue "E3" "  +0.88%  "

# Row 5 - BNB
Set-TextValue "D5" "694.89"
Set-TextValue "E5" "  +1.73%  "

# Row 6 - Solana
Set-TextValue "D6" "173.67"
Set-TextValue "E6" "  +1.96%  "

# Row 7 - LidoStakedEther
Set-TextValue "D7" "3.849.02"
Set-TextValue "E7" "  +0.93%  "

# Row 9 - XRP
Set-TextValue "E9" "  +0.09%  "

# Row 10 - Dogecoin
Set-TextValue "E10" "  +1.18%  "

# Row 11 - Toncoin
Set-TextValue "D11" "7.26"
Set-TextValue "E11" "  +1.14%  "

# Row 12 - Cardano
Set-TextValue "E12" "  -0.22%  "

# Row 13 - ShibaInu
Set-TextValue "E13" "  +4.60%  "

# Row 14 - Avalanche
Set-TextValue "D14" "36.41"
Set-TextValue "E14" "  +1.39%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "4.497.70"
Set-TextValue "E15" "  +0.88%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "3.845.20"
Set-TextValue "E16" "  +0.80%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "71.241.59"
Set-TextValue "E17" "  +0.84%  "

# Row 18 - Chainlink
Set-TextValue "E18" "  +0.03%  "

# Row 19 - Polkadot
Set-TextValue "E19" "  +0.50%  "

# Row 20 - TRON
Set-TextValue "D20" "0.114"
Set-TextValue "E20" "  +0.05%  "

# Row 21 - Uniswap
Set-TextValue "D21" "11.17"
Set-TextValue "E21" "  -0.45%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "494.03"
Set-TextValue "E22" "  +3.56%  "

# Row 23 - Polygon
Set-TextValue "E23" "  +1.10%  "

# Row 24 - Litecoin
Set-TextValue "D24" "85.08"
Set-TextValue "E24" "  +2.16%  "

# Row 25 - PEPE
Set-TextValue "E25" "  +1.54%  "

# Row 26 - InternetComputer(DFINITY)
Set-TextValue "D26" "12.31"
Set-TextValue "E26" "  +0.29%  "

# Row 27 - RenderToken
Set-TextValue "D27" "10.59"
Set-TextValue "E27" "  +2.43%  "

# Row 28 - Fetch.AI
Set-TextValue "E28" "  +1.81%  "

# Row 29 - WrappedeETH
Set-TextValue "D29" "4.003.82"
Set-TextValue "E29" "  +0.92%  "

# Row 30 - PancakeSwap
Set-TextValue "D30" "3.17"
Set-TextValue "E30" "  +8.34%  "

# Row 31 - Dai
Set-TextValue "E31" "  -0.09%  "

# Row 33 - ImmutableX
Set-TextValue "E33" "  -0.29%  "

# Row 34 - EthereumClassic
Set-TextValue "D34" "29.64"
Set-TextValue "E34" "  +0.25%  "

# Row 35 - Kaspa
Set-TextValue "D35" "0.180"
Set-TextValue "E35" "  -0.45%  "

# Row 36 - Aptos
Set-TextValue "D36" "9.27"
Set-TextValue "E36" "  +1.32%  "

# Row 37 - RenzoRestakedETH
Set-TextValue "D37" "3.800.89"
Set-TextValue "E37" "  +0.77%  "

# Row 38 - Binance-PegBSC-USD
Set-TextValue "D38" "1.00"
Set-TextValue "E38" "  -0.05%  "

# Row 39 - Hedera
Set-TextValue "D39" "0.104"
Set-TextValue "E39" "  +2.05%  "

# Row 40 - Stacks
Set-TextValue "D40" "2.39"
Set-TextValue "E40" "  +12.86%  "

# Row 41 - dogwifhat
Set-TextValue "E41" "  -0.04%  "

# Row 42 - Filecoin
Set-TextValue "E42" "  +1.74%  "

# Row 43 - Mantle
Set-TextValue "E43" "  +6.56%  "

# Row 44 - FirstDigitalUSD
Set-TextValue "E44" "  -0.10%  "

# Row 46 - Monero
Set-TextValue "D46" "163.37"
Set-TextValue "E46" "  +2.48%  "

# Row 47 - FLOKI
Set-TextValue "E47" "  +1.73%  "

# Row 48 - OKB
Set-TextValue "D48" "48.62"
Set-TextValue "E48" "  +0.90%  "

# Row 49 - Arweave
Set-TextValue "D49" "44.21"
Set-TextValue "E49" "  -3.87%  "

# Rows 50 & 51 swap order: Bittensor now ranks above TheGraph.
# Row 50 becomes Bittensor (was TheGraph)
Set-TextValue "B50" "Bittensor"
Set-TextValue "C50" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D50" "418.97"
Set-TextValue "E50" "  +5.40%  "

# Row 51 becomes TheGraph (was Bittensor)
Set-TextValue "B51" "TheGraph"
Set-TextValue "C51" "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D51" "0.303"
Set-TextValue "E51" "  +1.07%  "
